$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (shifts existing rows 11-20 down to 12-21)
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the weekly data
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "Macroferia Regional de Talca"
$ws.Range("C11").Value = "Maule"
$ws.Range("D11").Value = 45062
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100104
$ws.Range("H11").Value = "Frutos de pepita"
$ws.Range("I11").Value = 100104001
$ws.Range("J11").Value = "Granada"
$ws.Range("K11").Value = "Wonderfull"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("Q11").Value = '$/caja 15 kilos granel'
$ws.Range("R11").Value = 'Provincia de Curicó'
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 15

# Apply the date number format (matches the other rows' Fecha column) to the new row's date cell
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
